$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 142-143; this shifts the existing rows 142:242
# down to 144:244, preserving all of their original values/formatting.
$ws.Rows("142:143").Insert()

# New row 142 (Primera) - new weekly record inserted before the former row 142
$ws.Range("A142").Value = 1
$ws.Range("B142").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C142").Value = "Arica y Parinacota"
$ws.Range("D142").Value = 44596
$ws.Range("E142").Value = 15
$ws.Range("F142").Value = 100114014
$ws.Range("G142").Value = "Betarraga"
$ws.Range("H142").Value = "Sin especificar"
$ws.Range("I142").Value = "Primera"
$ws.Range("J142").Value = 900
$ws.Range("K142").Value = 400
$ws.Range("L142").Value = 450
$ws.Range("M142").Value = 425
$ws.Range("N142").Value = "`$/paquete 4 unidades"
$ws.Range("O142").Value = "Región de Arica y Parinacota"
$ws.Range("P142").Value = 106
$ws.Range("Q142").Value = 4
$ws.Range("R142").Value = "Hortaliza"

# New row 143 (Segunda) - companion record for the same new week
$ws.Range("A143").Value = 1
$ws.Range("B143").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C143").Value = "Arica y Parinacota"
$ws.Range("D143").Value = 44596
$ws.Range("E143").Value = 15
$ws.Range("F143").Value = 100114014
$ws.Range("G143").Value = "Betarraga"
$ws.Range("H143").Value = "Sin especificar"
$ws.Range("I143").Value = "Segunda"
$ws.Range("J143").Value = 1200
$ws.Range("K143").Value = 400
$ws.Range("L143").Value = 450
$ws.Range("M143").Value = 425
$ws.Range("N143").Value = "`$/paquete 5 unidades"
$ws.Range("O143").Value = "Región de Arica y Parinacota"
$ws.Range("P143").Value = 85
$ws.Range("Q143").Value = 5
$ws.Range("R143").Value = "Hortaliza"

# New rows 243-244 mirror what used to be the final rows (241-242, the
# last "Primera"/"Segunda" pair) which, after the shift, continue to
# repeat once more as the very last two rows of the table.
$ws.Range("A243").Value = 1
$ws.Range("B243").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C243").Value = "Arica y Parinacota"
$ws.Range("D243").Value = 44552
$ws.Range("E243").Value = 15
$ws.Range("F243").Value = 100114014
$ws.Range("G243").Value = "Betarraga"
$ws.Range("H243").Value = "Sin especificar"
$ws.Range("I243").Value = "Primera"
$ws.Range("J243").Value = 1200
$ws.Range("K243").Value = 300
$ws.Range("L243").Value = 350
$ws.Range("M243").Value = 325
$ws.Range("N243").Value = "`$/paquete 4 unidades"
$ws.Range("O243").Value = "Región de Arica y Parinacota"
$ws.Range("P243").Value = 81
$ws.Range("Q243").Value = 4
$ws.Range("R243").Value = "Hortaliza"

$ws.Range("A244").Value = 1
$ws.Range("B244").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C244").Value = "Arica y Parinacota"
$ws.Range("D244").Value = 44552
$ws.Range("E244").Value = 15
$ws.Range("F244").Value = 100114014
$ws.Range("G244").Value = "Betarraga"
$ws.Range("H244").Value = "Sin especificar"
$ws.Range("I244").Value = "Segunda"
$ws.Range("J244").Value = 1200
$ws.Range("K244").Value = 300
$ws.Range("L244").Value = 350
$ws.Range("M244").Value = 325
$ws.Range("N244").Value = "`$/paquete 5 unidades"
$ws.Range("O244").Value = "Región de Arica y Parinacota"
$ws.Range("P244").Value = 65
$ws.Range("Q244").Value = 5
$ws.Range("R244").Value = "Hortaliza"
